# Aggiornamento 15, 16, 17 marzo
# Adds three new daily rows (227-229) to the sheet, continuing the date
# series in column A (serial dates 44301, 44302, 44303) with zero values
# in columns B, C and D, matching the formatting of the preceding row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDates = @(44301, 44302, 44303)
$firstNewRow = 227
$templateRow = $firstNewRow - 1   # row 226, used as formatting template

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $firstNewRow + $i

    $ws.Cells.Item($row, 1).Value = $newDates[$i]
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0

    # Copy the style (s="2") of the date cell from the template row so the
    # new date cells are formatted the same way as the existing ones.
    $ws.Cells.Item($templateRow, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
}
